$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for rows 2-51 (columns B:E): 0-indexed rank in column A (unchanged),
# Coin name, Link, Price, Volume(1h). Values that would otherwise be
# auto-parsed as numbers by Excel are prefixed with a leading apostrophe so
# they stay text, matching the original inlineStr cells.
$data = @(
    @(0, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "60.639.37", "  -3.95%  "),
    @(1, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "2.921.16", "  -3.29%  "),
    @(2, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "'1.00", "  -0.04%  "),
    @(3, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "'528.88", "  -4.75%  "),
    @(4, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "'145.32", "  -6.38%  "),
    @(5, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "'1.00", "  +0.09%  "),
    @(6, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "'0.557", "  +0.08%  "),
    @(7, "LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "2.922.86", "  -3.59%  "),
    @(8, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "'0.109", "  -2.91%  "),
    @(9, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "'5.94", "  -6.99%  "),
    @(10, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "'0.355", "  -2.57%  "),
    @(11, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "3.434.25", "  -3.26%  "),
    @(12, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "'0.125", "  +1.73%  "),
    @(13, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "60.733.19", "  -3.90%  "),
    @(14, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "'23.00", "  -4.30%  "),
    @(15, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "2.929.00", "  -3.24%  "),
    @(16, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "'0.0000142", "  -5.40%  "),
    @(17, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "'5.00", "  -1.54%  "),
    @(18, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "'11.69", "  -2.41%  "),
    @(19, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "'364.52", "  -8.13%  "),
    @(20, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "'6.52", "  -2.05%  "),
    @(21, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "'0.999", "  -0.10%  "),
    @(22, "LEO", "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo", "'5.68", "  -2.08%  "),
    @(23, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "'64.34", "  -1.64%  "),
    @(24, "WrappedeETH", "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth", "3.062.98", "  -3.00%  "),
    @(25, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "'0.454", "  -2.00%  "),
    @(26, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "'0.184", "  -2.82%  "),
    @(27, "Binance-PegBSC-USD", "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd", "'0.999", "  +0.14%  "),
    @(28, "PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0₃0881", "  -9.89%  "),
    @(29, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "'7.83", "  -9.50%  "),
    @(30, "USDe", "https://coinranking.com/coin/exbfr2U-0+usde-usde", "'1.00", "  +0.06%  "),
    @(31, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "'1.67", "  -4.40%  "),
    @(32, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "'19.83", "  -2.57%  "),
    @(33, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "'157.27", "  -3.05%  "),
    @(34, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "'4.41", "  -6.60%  "),
    @(35, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "'5.66", "  -6.03%  "),
    @(36, "Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "'1.01", "  -9.91%  "),
    @(37, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "'1.22", "  -6.58%  "),
    @(38, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "'37.82", "  +0.11%  "),
    @(39, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "2.360.95", "  -6.68%  "),
    @(40, "Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "'1.48", "  -6.39%  "),
    @(41, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "'3.73", "  -4.72%  "),
    @(42, "Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "'0.647", "  -3.14%  "),
    @(43, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "'20.94", "  -7.77%  "),
    @(44, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "'0.0573", "  -4.47%  "),
    @(45, "FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "'0.998", "  -0.08%  "),
    @(46, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "'4.95", "  -2.55%  "),
    @(47, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "'0.0234", "  -5.41%  "),
    @(48, "WhiteBITCoin", "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt", "'10.40", "  -0.88%  "),
    @(49, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "'0.0930", "  -1.14%  ")
)

$rowIndex = 2
foreach ($r in $data) {
    $ws.Cells.Item($rowIndex, 2).Value = $r[1]
    $ws.Cells.Item($rowIndex, 3).Value = $r[2]
    $ws.Cells.Item($rowIndex, 4).Value = $r[3]
    $ws.Cells.Item($rowIndex, 5).Value = $r[4]
    $rowIndex++
}

Write-Host "Updated $($data.Count) rows"
